$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

# --- Company name / address changes: Sweet Entertainment Kft. -> 4 ÉP-SZAK 2000 Kft. ---

# Hozzájárulás paragraph: "a(z) Sweet Entertainment Kft. (, cégj.: ) székhelyéül"
Replace-Text "Sweet Entertainment Kft. (, cégj.: ) székhelyéül" "4 ÉP-SZAK 2000 Kft. (1139 Budapest, Országbíró u 2. 4. em. 19., cégj.: 01-09-687257) székhelyéül"

# "A Sweet Entertainment Kft. az ingatlan felett..."
Replace-Text "A Sweet Entertainment Kft. az ingatlan" "A 4 ÉP-SZAK 2000 Kft. az ingatlan"

# Bérlő description repeated in the two "amely létrejött egyrészről..." paragraphs
Replace-Text "Sweet Entertainment Kft. (, cégj.: ), képviseli: Albert János,  (, ig.sz.: , an.neve: 1), mint Bérlő" "4 ÉP-SZAK 2000 Kft. (1139 Budapest, Országbíró u 2. 4. em. 19., cégj.: 01-09-687257), képviseli: Bazsika István, ügyvezetõ (1012 Budapest, Logodi u. 48. fszt. 1., ig.sz.: 457361HA, an.neve: Süle Mária Margit), mint Bérlő"

# --- Dates: 2015-06-02 -> 2015-01-25 (covers every "Budapest, 2015-06-02" and "Az induló dátum: 2015-06-02") ---
Replace-Text "2015-06-02" "2015-01-25"

# Next due-date: 2015-12-02 -> 2015-07-25
Replace-Text "2015-12-02" "2015-07-25"

# --- Iratőrzési cím ---
Replace-Text "11. A Megbízott cégiratokat nem őriz. A cég iratainak őrzési címe: , " "11. A Megbízott cégiratokat nem őriz. A cég iratainak őrzési címe: Bazsika István, 1012 Budapest, Logodi u. 48. fszt. 1."

# --- Fees ---
Replace-Text "1. A megbízási/cégképviseleti díj nettó 0 Ft, azaz 1 Forint." "1. A megbízási/cégképviseleti díj nettó 32940 Ft, azaz Harminckettõezerkilencszáznegyven Forint."

Replace-Text "4. 1 (azaz egy) havi megbízási/cégképviseleti díj összege a szerződés aláírásakor nettó 0 Forint." "4. 1 (azaz egy) havi megbízási/cégképviseleti díj összege a szerződés aláírásakor nettó 5490 Forint."

# --- Contact details in annex ---
Replace-Text "Telefon: 06-30/546-34-12" "Telefon: 06-20/988-93-56"
Replace-Text "Email cím: sweet@infomail.com" "Email cím: bazsika.istvan@gmail.com"
Replace-Text "Postacím: , " "Postacím: Bazsika István, 1012 Budapest, Logodi u. 48. fszt. 1."
